# Add new columns I (I0) and J (IF) to the sheet, matching existing header style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, using same style as existing headers (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2..62 (I = I0 column, J = IF column)
$iValues = @(7,7,7,7,5,8,9,9,6,8,8,6,6,7,8,6,9,6,6,7,7,8,5,5,8,5,8,7,7,6,6,3,6,7,8,7,6,6,7,4,6,5,8,7,8,9,7,8,7,8,7,8,4,7,6,3,4,5,6,6,4)
$jValues = @(8,7,7,7,6,8,9,9,7,8,9,7,6,7,8,7,9,6,6,7,7,8,7,5,8,5,8,7,7,6,7,5,7,8,8,7,6,7,7,4,6,6,8,7,9,9,8,9,7,8,7,8,4,7,6,4,4,5,6,6,4)

for ($r = 2; $r -le 62; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
